$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert two new rows above the current row 28 ("Ca04a"/cd20), shifting the
# rest of the table down by two rows.
$ws.Rows.Item(28).Insert()
$ws.Rows.Item(28).Insert()

# Populate the two new rows with the new treatment variables, column by
# column (matches the shared-string insertion order of the source edit).
$ws.Range("A28").Value = "Ca04b"
$ws.Range("A29").Value = "Ca04c"

$ws.Range("B28").Value = "btki"
$ws.Range("B29").Value = "venet"

$ws.Range("C28").Value = "Cancer treatment"
$ws.Range("C29").Value = "Cancer treatment"

$ws.Range("D28").Value = "Most recent line of therapy includes BTK inhibitor"
$ws.Range("D29").Value = "Most recent line of therapy includes venetoclax"

$ws.Range("E28").Value = "0 = No; 1 = Yes"
$ws.Range("E29").Value = "0 = No; 1 = Yes"

# Grow the table (Table1) so it covers the two new rows.
$lo.Resize($ws.Range("A1:E233"))

# Update the saved view/selection to match the edited workbook.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A29").Select()
